$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new update date (this also auto-updates
# the workbook-level defined name "Lægesystemer___sende" to point at the
# renamed sheet).
$ws.Name = "Opdateret d. 05-12-2025"

# Vena (column I) has been approved ("Godkendt") for:
#  - Care Plan / CPD-DK (2.0)  -> row 8
#  - XDS Metadata              -> row 36
$ws.Range("I8").Value = "Godkendt"
$ws.Range("I36").Value = "Godkendt"
